$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 23: "Akurasi (%)" summary row, inserted above the (now blank) row 24
# and below the existing data table (rows 3-22).
# ---------------------------------------------------------------------------

# --- A23 (merged with B23): label cell, reuses the table's left-edge corner
#     border (same shape as C1) plus bold font and horizontal-only centering.
$ws.Range("C1").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").HorizontalAlignment = -4108
$ws.Range("A23").VerticalAlignment = -4107
$ws.Range("A23").Value = "Akurasi (%)"

# --- B23: empty cell, reuses the table's right-edge corner border (same
#     shape as V1) plus bold font and horizontal-only centering.
$ws.Range("V1").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Font.Bold = $true
$ws.Range("B23").HorizontalAlignment = -4108
$ws.Range("B23").VerticalAlignment = -4107

$ws.Range("A23:B23").Merge()

# --- C23: left=none / right=thin / top=medium / bottom=medium, bold, no
#     alignment override.
$ws.Range("A1").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Borders.Item(7).LineStyle = -4142
$cR = $ws.Range("C23").Borders.Item(10)
$cR.Weight = 2
$cR.ColorIndex = 1
$cB = $ws.Range("C23").Borders.Item(9)
$cB.Weight = -4138
$cB.ColorIndex = 1
$ws.Range("C23").Font.Bold = $true
$ws.Range("C23").Value = 20

# --- D23:U23: left=thin / right=thin / top=medium / bottom=medium, bold.
$midCols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")
$midVals = @(40,40,60,100,100,0,60,100,100,20,80,20,80,60,60,60,80,40)
for ($i = 0; $i -lt $midCols.Length; $i++) {
    $col = $midCols[$i]
    $cell = $col + "23"
    $ws.Range("A1").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $bL = $ws.Range($cell).Borders.Item(7)
    $bL.Weight = 2
    $bL.ColorIndex = 1
    $bR = $ws.Range($cell).Borders.Item(10)
    $bR.Weight = 2
    $bR.ColorIndex = 1
    $bB = $ws.Range($cell).Borders.Item(9)
    $bB.Weight = -4138
    $bB.ColorIndex = 1
    $ws.Range($cell).Font.Bold = $true
    $ws.Range($cell).Value = $midVals[$i]
}

# --- V23: left=thin / right=medium / top=medium / bottom=medium, bold.
$ws.Range("A1").Copy()
$ws.Range("V23").PasteSpecial(-4122)
$vL = $ws.Range("V23").Borders.Item(7)
$vL.Weight = 2
$vL.ColorIndex = 1
$vR = $ws.Range("V23").Borders.Item(10)
$vR.Weight = -4138
$vR.ColorIndex = 1
$vB = $ws.Range("V23").Borders.Item(9)
$vB.Weight = -4138
$vB.ColorIndex = 1
$ws.Range("V23").Font.Bold = $true
$ws.Range("V23").Value = 40

# ---------------------------------------------------------------------------
# New row 29: overall accuracy formula beneath the existing summary rows.
# ---------------------------------------------------------------------------
$ws.Range("D29").Formula = "=SUM(C23:V23)/20"
$ws.Range("E29").Value = "%"

# ---------------------------------------------------------------------------
# Misc view-state tweaks captured in the diff.
# ---------------------------------------------------------------------------
$ws.Range("K27").Select()
